$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is bumped by one day
# (46060 -> 46061) for every data row (rows 2 through 417).
$ws.Range("C2:C417").Value = 46061
